$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 79 (previously the "latest" row, formatted YYYY-MM-DD) reverts to the
# standard date format used by all the other data rows (YYYY-MM-DD HH:MM:SS)
$ws.Range("A79").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add new row 80 with the daily update values
$ws.Range("A80").Value = 45667
$ws.Range("B80").Value = 188
$ws.Range("C80").Value = 184
$ws.Range("D80").Value = 188

# Row 80's date cell becomes the new "latest row", taking on the special
# YYYY-MM-DD format that row 79 used to have
$ws.Range("A80").NumberFormat = "YYYY-MM-DD"
